# Update "想去人数" (F column) figures with freshly scraped totals.
# Output generated at 456a3b4

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 6970
$ws1.Range("F4").Value  = 58
$ws1.Range("F7").Value  = 6851
$ws1.Range("F8").Value  = 0
$ws1.Range("F9").Value  = 0
$ws1.Range("F10").Value = 0
$ws1.Range("F12").Value = 0
$ws1.Range("F13").Value = 408
$ws1.Range("F15").Value = 0
$ws1.Range("F16").Value = 0
$ws1.Range("F17").Value = 0
$ws1.Range("F18").Value = 40
$ws1.Range("F19").Value = 0
$ws1.Range("F20").Value = 5226
$ws1.Range("F23").Value = 0

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 6970
$ws4.Range("F3").Value  = 100
$ws4.Range("F4").Value  = 58
$ws4.Range("F5").Value  = 454
$ws4.Range("F7").Value  = 6851
$ws4.Range("F8").Value  = 74
$ws4.Range("F9").Value  = 0
$ws4.Range("F12").Value = 0
$ws4.Range("F13").Value = 408
$ws4.Range("F14").Value = 0
$ws4.Range("F17").Value = 0
$ws4.Range("F18").Value = 40
$ws4.Range("F19").Value = 0
$ws4.Range("F21").Value = 5226
$ws4.Range("F23").Value = 119
$ws4.Range("F24").Value = 166
$ws4.Range("F25").Value = 643
$ws4.Range("F26").Value = 0
$ws4.Range("F27").Value = 236
